# VL45 - Quantity Data Model
# 1) Duplicate the last ("Nädal") template sheet, placing the copy just
#    before it, and rename the copy to "Nädal9".
# 2) Fill in the new time-log entry on the freshly created "Nädal9" sheet.
# 3) Restore/adjust sheet view state (active sheet, scroll position, zoom).

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Nädal")

# Copy the template sheet, placing the new copy immediately before itself.
# (The original "Nädal" sheet is pushed one slot later; it keeps its name
# and stays an empty template for the next week.)
$template.Copy($template)

# The newly inserted copy lands right before the original "Nädal" sheet.
$newSheet = $wb.Worksheets.Item("Nädal (2)")
$newSheet.Name = "Nädal9"

# Fill in the actual log entry on the newly created "Nädal9" sheet.
$newSheet.Range("B7").Value = (Get-Date -Year 2020 -Month 3 -Day 26 -Hour 0 -Minute 0 -Second 0)
$newSheet.Range("C7").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 11 -Minute 38 -Second 0)
$newSheet.Range("D7").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 12 -Minute 11 -Second 0)
$newSheet.Range("G7").Value = "VL45"
$newSheet.Range("H7").Value = "Quantity Data Model"

# Make "Nädal9" the active sheet with the expected view state.
$newSheet.Activate()
$newSheet.Application.ActiveWindow.Zoom = 87
$newSheet.Range("E7").Select()

# "Nädal8" view loses its selected/top-left state.
$week8 = $wb.Worksheets.Item("Nädal8")
$week8.Activate()
$week8.Application.ActiveWindow.ScrollRow = 5

# Finally land back on "Nädal9" to match the saved active tab.
$newSheet.Activate()
